$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (and fix the PEPE/Litecoin row
# ordering in rows 24-25) to match the latest scrape.
# Price (column D) values are entered with a leading apostrophe so Excel keeps
# them as literal text (matching the original inlineStr cells) instead of
# auto-converting number-like strings (e.g. "706.69") into floating point values.
$ws.Range("D2").Value = "'71.256.13"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").Value = "'3.814.33"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'706.69"
$ws.Range("E5").Value = "  +12.35%  "
$ws.Range("D6").Value = "'173.91"
$ws.Range("E6").Value = "  +5.15%  "
$ws.Range("D7").Value = "'3.813.37"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("E11").Value = "  +10.07%  "
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +9.65%  "
$ws.Range("D14").Value = "'36.39"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").Value = "'4.457.70"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "'3.812.90"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "'71.241.80"
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'7.28"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "'11.23"
$ws.Range("E21").Value = "  +18.10%  "
$ws.Range("D22").Value = "'484.36"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000149"
$ws.Range("E24").Value = "  +6.08%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'84.01"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("D27").Value = "'10.66"
$ws.Range("E27").Value = "  +5.52%  "
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").Value = "'3.966.77"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'3.11"
$ws.Range("E31").Value = "  +16.85%  "
$ws.Range("D32").Value = "'2.30"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "'7.57"
$ws.Range("E33").Value = "  +6.80%  "
$ws.Range("D34").Value = "'29.66"
$ws.Range("E34").Value = "  +4.39%  "
$ws.Range("D35").Value = "'0.179"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").Value = "'9.27"
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'3.765.64"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("E40").Value = "  +8.62%  "
$ws.Range("D41").Value = "'6.00"
$ws.Range("E41").Value = "  +4.04%  "
$ws.Range("D42").Value = "'0.000341"
$ws.Range("E42").Value = "  +31.92%  "
$ws.Range("D43").Value = "'2.24"
$ws.Range("E43").Value = "  +13.07%  "
$ws.Range("D44").Value = "'0.975"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'161.16"
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("D48").Value = "'45.33"
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("D49").Value = "'49.35"
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'0.301"
$ws.Range("E51").Value = "  +2.97%  "
